$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are stored as text in the source sheet, so force
# text entry with a leading apostrophe to avoid Excel auto-converting them
# to numbers (which would drop meaningful trailing zeros). Reset the style
# afterwards so the forced "quote prefix" formatting doesn't linger.
$ws.Range("D3").Value  = "'22.93"
$ws.Range("D3").Style  = "Normal"

$ws.Range("D5").Value  = "'0.05643"
$ws.Range("D5").Style  = "Normal"

$ws.Range("D6").Value  = "'3.425"
$ws.Range("D6").Style  = "Normal"

$ws.Range("D9").Value  = "'0.9179"
$ws.Range("D9").Style  = "Normal"

$ws.Range("D10").Value = "'0.01156"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = "9OneONEBestin24h"

$ws.Range("D11").Value = "'0.1437"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").Value = "'0.07510"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").Value = "'0.03125"
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").Value = "'0.03110"
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").Value = "'0.09352"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").Value = "'3.557"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").Value = "'0.001582"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").Value = "'0.04768"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").Value = "'0.006366"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").Value = "'0.005001"
$ws.Range("D20").Style = "Normal"

$ws.Range("D24").Value = "'2.191"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").Value = "'0.3299"
$ws.Range("D25").Style = "Normal"

$ws.Range("D28").Value = "'0.0003031"
$ws.Range("D28").Style = "Normal"

$ws.Range("D40").Value = "'0.04042"
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").Value = "'0.006797"
$ws.Range("D41").Style = "Normal"

$ws.Range("D43").Value = "'0.002709"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").Value = "'0.007551"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").Value = "'0.00005801"
$ws.Range("D45").Style = "Normal"

$ws.Range("D47").Value = "'0.4998"
$ws.Range("D47").Style = "Normal"

$ws.Range("E48").Value = "47BOLOBOLO"

$ws.Range("D50").Value = "'0.01010"
$ws.Range("D50").Style = "Normal"
